# Applies the cryptos-list price/volume update described in the commit
# "Updated cryptos list on Sat Jan  6 06:19:31 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds dotted numeric-looking text (e.g. "43.906.31", "1.00").
# Force it to Text format first so Excel does not silently reinterpret/round it
# as a number when we assign the literal string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.906.31"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.233.87"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.19"
$ws.Range("E5").Value = "  -4.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.15"
$ws.Range("E6").Value = "  -7.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("E7").Value = "  -1.53%  "

$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -6.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.09"
$ws.Range("E10").Value = "  -7.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  -4.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.08"
$ws.Range("E12").Value = "  -8.24%  "

$ws.Range("E13").Value = "  -2.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.574.13"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.271.52"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("E16").Value = "  -5.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.47"
$ws.Range("E17").Value = "  -5.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.714.38"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0953"
$ws.Range("E19").Value = "  -3.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.24"
$ws.Range("E20").Value = "  -8.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.14"
$ws.Range("E21").Value = "  -6.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.34"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.54"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("E24").Value = "  -7.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.77"
$ws.Range("E27").Value = "  -3.53%  "

$ws.Range("E28").Value = "  -3.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.10"
$ws.Range("E29").Value = "  -3.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.91"
$ws.Range("E30").Value = "  -5.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.78"
$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.80"
$ws.Range("E32").Value = "  -4.71%  "

$ws.Range("E33").Value = "  -6.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.25"
$ws.Range("E34").Value = "  +4.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.62"
$ws.Range("E35").Value = "  -2.87%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("E37").Value = "  -9.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.72"
$ws.Range("E38").Value = "  -11.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.38"
$ws.Range("E39").Value = "  -8.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.77"
$ws.Range("E40").Value = "  -10.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0297"
$ws.Range("E41").Value = "  -6.04%  "

$ws.Range("E42").Value = "  -13.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.739.78"
$ws.Range("E44").Value = "  -3.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.81"
$ws.Range("E45").Value = "  +1.35%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.185"
$ws.Range("E46").Value = "  -6.92%  "

$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.94"
$ws.Range("E47").Value = "  -5.05%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "14.84"
$ws.Range("E48").Value = "  +5.72%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.36"
$ws.Range("E49").Value = "  -4.18%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.07"
$ws.Range("E50").Value = "  -3.90%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.68"
$ws.Range("E51").Value = "  -10.63%  "
